$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("B3").Value = "Name"
$ws.Range("E4:E7").Style = "Normal"
$ws.Range("B3").Select() | Out-Null
